$p = $ppt.ActivePresentation

# --- 1) Update the cached Date / Slide-number placeholder text on every
#        slide layout served by the slide master (date field + slide-number
#        field are present on all 11 layouts in this template). ---
$newDate = "05.01.2025"
$newSlideNum = [string]([char]0x2039) + "Nr." + [string]([char]0x203A)

$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $lay = $layouts.Item($li)
    $shapes = $lay.Shapes
    for ($si = 1; $si -le $shapes.Count; $si++) {
        $shp = $shapes.Item($si)
        if ($shp.HasTextFrame -and $shp.PlaceholderFormat.Type -eq 16) {
            # ppPlaceholderDate
            $shp.TextFrame.TextRange.Text = $newDate
        }
        elseif ($shp.HasTextFrame -and $shp.PlaceholderFormat.Type -eq 13) {
            # ppPlaceholderSlideNumber
            $shp.TextFrame.TextRange.Text = $newSlideNum
        }
    }
}

# --- 2) Rename "Content Placeholder 2" (the free-floating textbox, not the
#        real placeholder, on slide 4) to "TestShape" and replace its body
#        text with two separate runs: "TestShape" + " " + "PlaceHolder". ---
$slide = $p.Slides.Item(4)
for ($si = 1; $si -le $slide.Shapes.Count; $si++) {
    $shp = $slide.Shapes.Item($si)
    if ($shp.Name -eq "Content Placeholder 2" -and $shp.Type -eq 17) {
        # msoTextBox (id 4) -- the plain textbox, not the real placeholder.
        $shp.Name = "TestShape"
        $tr = $shp.TextFrame.TextRange
        $tr.Text = "TestShape PlaceHolder"
    }
}
